$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 17, 18, 19: cyclic rotation ---
# new17 = old18, new18 = old19, new19 = old17
$row17 = $ws.Range("F17:V17").Value2
$row18 = $ws.Range("F18:V18").Value2
$row19 = $ws.Range("F19:V19").Value2

$ws.Range("F17:V17").Value = $row18
$ws.Range("F18:V18").Value = $row19
$ws.Range("F19:V19").Value = $row17

# --- Rows 20, 21: swap ---
$row20 = $ws.Range("F20:V20").Value2
$row21 = $ws.Range("F21:V21").Value2

$ws.Range("F20:V20").Value = $row21
$ws.Range("F21:V21").Value = $row20

# --- Rows 44, 45: swap ---
$row44 = $ws.Range("F44:V44").Value2
$row45 = $ws.Range("F45:V45").Value2

$ws.Range("F44:V44").Value = $row45
$ws.Range("F45:V45").Value = $row44

# --- Rows 67, 68: swap ---
$row67 = $ws.Range("F67:V67").Value2
$row68 = $ws.Range("F68:V68").Value2

$ws.Range("F67:V67").Value = $row68
$ws.Range("F68:V68").Value = $row67

# --- Rows 77, 78: swap ---
$row77 = $ws.Range("F77:V77").Value2
$row78 = $ws.Range("F78:V78").Value2

$ws.Range("F77:V77").Value = $row78
$ws.Range("F78:V78").Value = $row77

# --- New row 154: append new match record ---
$ws.Range("A154").Value = 153
$ws.Range("B154").Value = "colombia"
$ws.Range("C154").Value = "primera-b"
# "2023" is a purely-numeric string; force text storage like the rest of column D
$ws.Range("D154").NumberFormat = "@"
$ws.Range("D154").Value = "2023"
$ws.Range("E154").Value = 45245.0625
$ws.Range("F154").Value = "Cucuta"
$ws.Range("G154").Value = 1
$ws.Range("H154").Value = "Fortaleza"
$ws.Range("I154").Value = 0
$ws.Range("J154").Value = 2.11
$ws.Range("K154").Value = "14/11/2023 13:42"
$ws.Range("L154").Value = 2.23
$ws.Range("M154").Value = "15/11/2023 01:29"
$ws.Range("N154").Value = 3.12
$ws.Range("O154").Value = "14/11/2023 13:42"
$ws.Range("P154").Value = 3.05
$ws.Range("Q154").Value = "15/11/2023 01:29"
$ws.Range("R154").Value = 3.84
$ws.Range("S154").Value = "14/11/2023 13:42"
$ws.Range("T154").Value = 3.72
$ws.Range("U154").Value = "15/11/2023 01:29"
$ws.Range("V154").Value = "https://www.betexplorer.com/football/colombia/primera-b/cucuta-fortaleza-c-e-i-f/f5tE3iBG/"

# Match the exact cell formatting used by the equivalent cells one row up
# (copy-format only: doesn't touch the values/types already written above,
# and also clears D154's temporary NumberFormat="@" override back down to
# the plain/default style while keeping its value stored as text)
$ws.Range("A153").Copy()
$ws.Range("A154").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D153").Copy()
$ws.Range("D154").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E153").Copy()
$ws.Range("E154").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
